$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.064.60"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'2.413.85"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'554.42"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'136.65"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'5.67"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "'24.81"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'2.845.73"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'59.984.62"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'2.415.39"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  +3.64%  "
$ws.Range("D20").Value = "'326.52"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'64.75"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +5.74%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +5.14%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'1.78"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "'170.89"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("E32").Value = "  +5.57%  "
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").Value = "'18.42"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "'324.31"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'146.45"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").Value = "'3.62"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'19.79"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("D45").Value = "'0.0516"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.938"
$ws.Range("E51").Value = "  -1.75%  "
